# Interdiff v1 -> v2: both "edit 5 …" rectangles become "clear".
#
#   - "Rectangle 20" (shape id 21)  -> text becomes "clear"
#   - "Rectangle 28" (shape id 29)  -> text becomes "clear"
#
# Slide 1 is the only slide in the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName {
    param($slide, $name)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

$rect20 = Get-ShapeByName $s "Rectangle 20"
$rect28 = Get-ShapeByName $s "Rectangle 28"

$rect20.TextFrame.TextRange.Text = "clear"
$rect28.TextFrame.TextRange.Text = "clear"
